$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6656782362548082
$ws.Range("C2").Value = 0.7098241526097382
$ws.Range("D2").Value = 1.356451150188942
$ws.Range("E2").Value = 0.6948499821685598
$ws.Range("F2").Value = 0.6767035828922799
$ws.Range("G2").Value = 0.739899615432439
$ws.Range("H2").Value = 0.6952771648018571
$ws.Range("B3").Value = 0.6945291719047449
$ws.Range("C3").Value = 0.7155518372566629
$ws.Range("D3").Value = 1.276233810667159
$ws.Range("E3").Value = 0.7010362124351743
$ws.Range("F3").Value = 0.6980272502158601
$ws.Range("G3").Value = 0.7431061246380334
$ws.Range("H3").Value = 0.7012035129671159
$ws.Range("B4").Value = 0.6115471694156789
$ws.Range("C4").Value = 0.7396515998568782
$ws.Range("D4").Value = 0.9202580287282726
$ws.Range("E4").Value = 0.6867376105788219
$ws.Range("F4").Value = 0.6244327983843749
$ws.Range("G4").Value = 0.7215458874735288
$ws.Range("H4").Value = 0.6869053491050978
$ws.Range("B5").Value = 0.7293867604872588
$ws.Range("C5").Value = 0.7668410534285477
$ws.Range("D5").Value = 0.4247129094724361
$ws.Range("E5").Value = 0.7046736107485284
$ws.Range("F5").Value = 0.7258354903946782
$ws.Range("G5").Value = 0.7436580617940645
$ws.Range("H5").Value = 0.7044098395282327
$ws.Range("B6").Value = 0.7187228668629096
$ws.Range("C6").Value = 0.7644877422023012
$ws.Range("D6").Value = 0.2837220892666605
$ws.Range("E6").Value = 0.7027173425921137
$ws.Range("F6").Value = 0.7154820880022621
$ws.Range("G6").Value = 0.7357746999847482
$ws.Range("H6").Value = 0.7024036911943025
$ws.Range("B7").Value = 0.7498245758194799
$ws.Range("C7").Value = 0.7856843859769028
$ws.Range("D7").Value = 0.401576885974789
$ws.Range("E7").Value = 0.7130282815819784
$ws.Range("F7").Value = 0.7460343174914137
$ws.Range("G7").Value = 0.7585522726615487
$ws.Range("H7").Value = 0.7127657133716443
$ws.Range("B8").Value = 0.6903222019640344
$ws.Range("C8").Value = 0.5948629548070007
$ws.Range("D8").Value = 0.3038021567545686
$ws.Range("E8").Value = 0.6848147166973333
$ws.Range("F8").Value = 0.686666898898007
$ws.Range("G8").Value = 0.6022878936784698
$ws.Range("H8").Value = 0.6845477254926853
$ws.Range("B9").Value = 0.7568206031162165
$ws.Range("C9").Value = 0.7453565595481404
$ws.Range("D9").Value = 0.5641100114228719
$ws.Range("E9").Value = 0.7103386006351842
$ws.Range("F9").Value = 0.7531255498779952
$ws.Range("G9").Value = 0.7261943490557077
$ws.Range("H9").Value = 0.7101543221519784
$ws.Range("B10").Value = 0.7523855043214686
$ws.Range("C10").Value = 0.7814954067393087
$ws.Range("D10").Value = 1.011569275504135
$ws.Range("E10").Value = 0.7119964603128214
$ws.Range("F10").Value = 0.7535477510794063
$ws.Range("G10").Value = 0.7743276069489899
$ws.Range("H10").Value = 0.7122758823361749
$ws.Range("B11").Value = 0.7447263700905355
$ws.Range("C11").Value = 0.7828152745360828
$ws.Range("D11").Value = 0.9020829115585566
$ws.Range("E11").Value = 0.7126442726389244
$ws.Range("F11").Value = 0.7456990317317909
$ws.Range("G11").Value = 0.7720118857993419
$ws.Range("H11").Value = 0.712790944671245
$ws.Range("B12").Value = 0.651220674327408
$ws.Range("C12").Value = 0.7856065710047644
$ws.Range("D12").Value = 0.391530093490966
$ws.Range("E12").Value = 0.7063020607180248
$ws.Range("F12").Value = 0.6540875005989942
$ws.Range("G12").Value = 0.7574339153657651
$ws.Range("H12").Value = 0.7060317019685548
$ws.Range("B13").Value = 0.7479563703930113
$ws.Range("C13").Value = 0.7805575362458407
$ws.Range("D13").Value = 0.8449974102226457
$ws.Range("E13").Value = 0.7111449616502787
$ws.Range("F13").Value = 0.7481255078652088
$ws.Range("G13").Value = 0.7686031626159526
$ws.Range("H13").Value = 0.7113333033618732
